# "All query requirements have been completed"
#
# The grading sheet has a "Grade" column (C) that gets the letter S/G/B/C
# filled in once a requirement's query has been verified. Four rows were
# still blank; now that those queries are done, fill them in with the
# appropriate grade letter, matching the formatting used by every other
# already-graded cell in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from an already-graded cell (C8) onto the blank ones so
# the new entries look like the rest of the column.
$ws.Range("C8").Copy()
$ws.Range("C18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C24").PasteSpecial(-4122)  # xlPasteFormats

# "1 Delete" / "1 Update" -> grade G
$ws.Range("C18").Value = "G"
$ws.Range("C19").Value = "G"

# "1 single-row" / "1 multiple-row" -> grade C
$ws.Range("C23").Value = "C"
$ws.Range("C24").Value = "C"

# Leave the sheet focused where the author last worked.
$ws.Range("C9").Select()
$excel.ActiveWindow.Zoom = 112
